# Update countries & provincias Spain
# Applies the data refresh described in the commit to the "Pais" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp banner (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 12:15"

# --- Rumania now overtakes Catar: new Rumania data moves into row 32, ---
# --- the old (previous) Catar values slide down into row 33.         ---
$ws.Range("A32").Value = "Rumania"
$ws.Range("B32").Value = 127572
$ws.Range("C32").Value = 2158
$ws.Range("D32").Value = 102476
$ws.Range("E32").Value = 20271
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 33
$ws.Range("H32").Value = 4825

$ws.Range("A33").Value = "Catar"
$ws.Range("B33").Value = 125533
$ws.Range("C33").Value = 0
$ws.Range("D33").Value = 122448
$ws.Range("E33").Value = 2871
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 214

# --- Polonia now overtakes Guatemala: new Polonia data moves into row 45, ---
# --- the old (previous) Guatemala values slide down into row 46.         ---
$ws.Range("A45").Value = "Polonia"
$ws.Range("B45").Value = 91514
$ws.Range("C45").Value = 1552
$ws.Range("D45").Value = 69695
$ws.Range("E45").Value = 19306
$ws.Range("F45").Value = 0
$ws.Range("G45").Value = 30
$ws.Range("H45").Value = 2513

$ws.Range("A46").Value = "Guatemala"
$ws.Range("B46").Value = 90968
$ws.Range("C46").Value = 0
$ws.Range("D46").Value = 79654
$ws.Range("E46").Value = 8076
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 3238

# --- Finlandia now overtakes Guayana Francesa: new Finlandia data moves ---
# --- into row 105, the old Guayana Francesa values slide into row 106. ---
$ws.Range("A105").Value = "Finlandia"
$ws.Range("B105").Value = 9992
$ws.Range("C105").Value = 100
$ws.Range("D105").Value = 7850
$ws.Range("E105").Value = 1797
$ws.Range("F105").Value = 0
$ws.Range("G105").Value = 0
$ws.Range("H105").Value = 345

$ws.Range("A106").Value = "Guayana Francesa"
$ws.Range("B106").Value = 9929
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 9569
$ws.Range("E106").Value = 294
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 0
$ws.Range("H106").Value = 66

# --- Standalone data refreshes (country/order unchanged) ---

# Alemania (row 25)
$ws.Range("D25").Value = 256000
$ws.Range("E25").Value = 24910

# Austria (row 67)
$ws.Range("B67").Value = 44813
$ws.Range("C67").Value = 772
$ws.Range("D67").Value = 35644
$ws.Range("E67").Value = 8370
$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 799

# Georgia (row 116)
$ws.Range("E116").Value = 3034
$ws.Range("G116").Value = 2
$ws.Range("H116").Value = 38

# Hong Kong (row 124)
$ws.Range("B124").Value = 5088
$ws.Range("C124").Value = 8
$ws.Range("D124").Value = 4827
$ws.Range("E124").Value = 156

# Gibraltar (row 180)
$ws.Range("B180").Value = 396
$ws.Range("C180").Value = 5
$ws.Range("D180").Value = 344
$ws.Range("E180").Value = 52
